$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.165764093399048
$ws.Range("B1").Value = 2.426128625869751
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.373390913009644
$ws.Range("E1").Value = 1.235076546669006
